$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row66 = @(576.20799999999997, 1179.5839999999998, 435.02899999999994, 650.00299999999993, 473.85599999999994, 820.72500000000002, 488.029, 240.227, 133.07599999999999, 223.92599999999999, 206.53300000000002, 187.64699999999999, 860.8950000000001, 966.72400000000016, 323.42900000000003, 457.43, 346.26499999999999, 250.32299999999998, 146.06, 54.919999999999995, 212.50200000000001, 82.54, 48.12, 18.740000000000002)

$row67 = @(454.21100000000001, 588.02700000000004, 190.36699999999999, 117.895, 18.335999999999999, 258.53199999999998, 176.58300000000003, 13.478, 166.113, 172.57699999999997, 85.188000000000002, 144.30000000000001, 177.34100000000001, 204.85199999999998, 209.23500000000001, 78.531000000000006, 59.597000000000001, 52.978000000000002, 97.835000000000008, 31.820999999999998, 0, 49.234000000000002, 12.827999999999999, 37.262)

$ws.Cells.Item(66, 1).Value = 45933
$ws.Cells.Item(66, 2).Value = "四方坪站充电量(kw)"
for ($i = 0; $i -lt $row66.Length; $i++) {
    $ws.Cells.Item(66, 3 + $i).Value = $row66[$i]
}

$ws.Cells.Item(67, 1).Value = 45933
$ws.Cells.Item(67, 2).Value = "高岭站充电量(kw)"
for ($i = 0; $i -lt $row67.Length; $i++) {
    $ws.Cells.Item(67, 3 + $i).Value = $row67[$i]
}

$ws.Range("A66:A67").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("C66:Z67").NumberFormat = "0.00_);[Red]\(0.00\)"

[void]$ws.Range("E68").Select()
